# ValueSet-KLInterventions: bump published term version 1.0.0 -> 1.1.0
# and refresh the publish Date to match (commit: "Added 1.1.0 of term").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: "Version" | "1.0.0" -> "1.1.0"
$ws.Range("B3").Value = "1.1.0"

# Row 8: "Date" | old timestamp -> new publish timestamp
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
